$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.418.92'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.929.14'
$ws.Range('E3').Value = '  +4.16%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9983'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('E5').Value = '  +3.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4762'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '44.34'
$ws.Range('E8').Value = '  +2.37%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2874'
$ws.Range('E9').Value = '  +4.45%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06582'
$ws.Range('E10').Value = '  +4.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.08'
$ws.Range('E11').Value = '  +8.04%  '
$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '107.69'
$ws.Range('E12').Value = '  +27.06%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.923.99'
$ws.Range('E13').Value = '  +3.81%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.07621'
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.151'
$ws.Range('E15').Value = '  +4.01%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.6600'
$ws.Range('E16').Value = '  +5.51%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '308.75'
$ws.Range('E17').Value = '  +25.36%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '30.427.85'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.02'
$ws.Range('E19').Value = '  +2.70%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9990'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.000007526'
$ws.Range('E21').Value = '  +2.70%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.168.96'
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.314'
$ws.Range('E23').Value = '  +8.29%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.9976'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '6.295'
$ws.Range('E25').Value = '  +6.66%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '167.83'
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.252'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.24'
$ws.Range('E28').Value = '  +12.58%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.037'
$ws.Range('E29').Value = '  +8.71%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.1117'
$ws.Range('E30').Value = '  +8.38%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.359'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.101'
$ws.Range('E32').Value = '  +1.46%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.933'
$ws.Range('E33').Value = '  +2.83%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05011'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7448'
$ws.Range('E35').Value = '  +6.55%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.151'
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.752'
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01951'
$ws.Range('E38').Value = '  +2.75%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.707'
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.057'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8790'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '107.14'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '70.56'
$ws.Range('E43').Value = '  +11.27%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.815'
$ws.Range('E44').Value = '  +5.58%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9988'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4175'
$ws.Range('E46').Value = '  +2.81%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.296'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.274'
$ws.Range('E48').Value = '  +8.02%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '34.90'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.1205'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05621'
$ws.Range('E51').Value = '  +2.13%  '
